$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.723.17'
$ws.Range('E2').Value = '  -1.52%  '
$ws.Range('D3').Value = '1.547.07'
$ws.Range('E3').Value = '  -1.67%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '206.02'
$ws.Range('E5').Value = '  -0.55%  '
$ws.Range('E6').Value = '  -1.87%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = '21.40'
$ws.Range('E8').Value = '  -3.88%  '
$ws.Range('E9').Value = '  -1.51%  '
$ws.Range('E10').Value = '  -1.21%  '
$ws.Range('D11').Value = '0.0852'
$ws.Range('E11').Value = '  -1.89%  '
$ws.Range('D12').Value = '1.767.84'
$ws.Range('E12').Value = '  -1.64%  '
$ws.Range('D13').Value = '1.548.17'
$ws.Range('E13').Value = '  -1.59%  '
$ws.Range('E14').Value = '  -2.79%  '
$ws.Range('E15').Value = '  -1.21%  '
$ws.Range('D16').Value = '26.730.80'
$ws.Range('D17').Value = '61.16'
$ws.Range('E17').Value = '  -1.79%  '
$ws.Range('D18').Value = '212.43'
$ws.Range('E18').Value = '  -0.87%  '
$ws.Range('E19').Value = '  +0.68%  '
$ws.Range('E20').Value = '  -1.86%  '
$ws.Range('E21').Value = '  -0.09%  '
$ws.Range('D22').Value = '4.07'
$ws.Range('E22').Value = '  -1.30%  '
$ws.Range('E23').Value = '  -5.73%  '
$ws.Range('E24').Value = '  -1.80%  '
$ws.Range('D25').Value = '152.82'
$ws.Range('E25').Value = '  +0.21%  '
$ws.Range('D26').Value = '6.51'
$ws.Range('E26').Value = '  -3.00%  '
$ws.Range('D27').Value = '14.85'
$ws.Range('E27').Value = '  -0.61%  '
$ws.Range('E29').Value = '  -1.78%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').Value = '1.10'
$ws.Range('E30').Value = '  -1.56%  '
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').Value = '0.0459'
$ws.Range('E31').Value = '  -0.90%  '
$ws.Range('E32').Value = '  +0.17%  '
$ws.Range('D33').Value = '1.338.25'
$ws.Range('E33').Value = '  -4.15%  '
$ws.Range('E34').Value = '  -0.43%  '
$ws.Range('E35').Value = '  -3.16%  '
$ws.Range('E36').Value = '  -0.66%  '
$ws.Range('D37').Value = '0.927'
$ws.Range('E37').Value = '  -1.52%  '
$ws.Range('E38').Value = '  -0.55%  '
$ws.Range('E39').Value = '  +0.95%  '
$ws.Range('D40').Value = '5.77'
$ws.Range('E40').Value = '  +6.42%  '
$ws.Range('E41').Value = '  -1.82%  '
$ws.Range('D42').Value = '0.995'
$ws.Range('E42').Value = '  -1.45%  '
$ws.Range('E43').Value = '  -0.14%  '
$ws.Range('E44').Value = '  -3.56%  '
$ws.Range('D45').Value = '62.70'
$ws.Range('E45').Value = '  -1.66%  '
$ws.Range('D46').Value = '1.681.99'
$ws.Range('E46').Value = '  -1.64%  '
$ws.Range('D47').Value = '2.25'
$ws.Range('E47').Value = '  -3.69%  '
$ws.Range('D48').Value = '85.88'
$ws.Range('E48').Value = '  +0.34%  '
$ws.Range('D49').Value = '0.0507'
$ws.Range('E49').Value = '  +2.69%  '
$ws.Range('D50').Value = '0.0₇0975'
$ws.Range('E50').Value = '  -0.49%  '
$ws.Range('D51').Value = '0.0953'
$ws.Range('E51').Value = '  +0.14%  '
